$wb = $excel.ActiveWorkbook

# --- Add the new "Scatter" worksheet, moved after the last existing sheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Scatter"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-fetch by name: after Move() the earlier handle tracks whatever sheet
# now sits in the original slot, not the "Scatter" sheet itself.
$ws = $wb.Worksheets.Item("Scatter")

# --- Fill the data: column A = 0..9, column B = squares of column A ---
$yValues = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9)
$xValues = @(0, 1, 4, 9, 16, 25, 36, 49, 64, 81)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $yValues[$i]
    $ws.Cells.Item($row, 2).Value = $xValues[$i]
}

# --- Add the scatter (XY) chart: lines + markers style, like the workbook's others ---
$ws.Shapes.AddChart2(-1, 74) | Out-Null
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart

$ser = $chart.SeriesCollection(1)
$ser.XValues = "='Scatter'!`$B`$1:`$B`$10"
$ser.Values = "='Scatter'!`$A`$1:`$A`$10"

# Axis 1 (document order) is the x-value axis (data range 0-90 for the
# squares in column B); axis 2 is the y-value axis (data range 0-9, column A).
$xAxis = $chart.Axes(1)
$xAxis.MinimumScale = 0
$xAxis.MaximumScale = 90
$xAxis.MajorUnit = 10
$xAxis.TickLabels.NumberFormat = "General"
$xAxis.TickLabels.NumberFormatLinked = 1

$yAxis = $chart.Axes(2)
$yAxis.MinimumScale = 0
$yAxis.MaximumScale = 10

Write-Output "Scatter sheet + chart added"
